$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("AF4").Value = 0.65
$ws.Range("AF5").Value = 0.8
$ws.Range("AF6").Value = 0.717
$ws.Range("AF7").Value = 0.765
$ws.Range("AF8").Value = 0.8
$ws.Range("AF9").Value = 0.8
$ws.Range("AF10").Value = 0.8
$ws.Range("AF11").Value = 0.8
$ws.Range("AF12").Value = 1
$ws.Range("AF13").Value = 1.8
